$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertTo-Excel2DArray($jagged) {
    $nrows = $jagged.Count
    $ncols = $jagged[0].Count
    $arr = New-Object 'object[,]' $nrows, $ncols
    for ($i = 0; $i -lt $nrows; $i++) {
        for ($j = 0; $j -lt $ncols; $j++) {
            $arr[$i, $j] = $jagged[$i][$j]
        }
    }
    return $arr
}

# Updated Ligand-side stats (columns G:J, rows 2-17) recomputed with the
# new TPM values; these depend only on each row's Sending cluster.
$gjJagged = @(
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981)
)

# Updated Receptor-side stats + Edge weight/specificity columns
# (columns M:T, rows 2-17) recomputed with the new TPM values.
$mtJagged = @(
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862,0.9637813576934446,8.674032219241001,0.002198043187798391,0.002335635670569921),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691,5.060714163080223,45.546427467722,0.01154169273223428,0.012264176333662),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412,1.844095751143833,11.064574506863,0.004205727856316295,0.002979331208913716),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981,12.68844067872533,114.195966108528,0.02893782949319845,0.03074927151949298),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862,5.060714163080223,45.546427467722,0.01154169273223428,0.012264176333662),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691,26.57327581194711,239.159482307524,0.06060421008321335,0.06439789520188881),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412,9.683152108507665,58.09891265104599,0.02208383297583021,0.01564415365072081),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981,66.62566248833066,599.6309623949759,0.1519494877090276,0.1614611785555975),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862,1.844095751143833,11.064574506863,0.004205727856316295,0.002979331208913716),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691,9.683152108507665,58.09891265104599,0.02208383297583021,0.01564415365072081),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412,3.528486115902249,14.113944463609,0.008047224412870478,0.003800427679819311),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981,24.278016333784,145.668098002704,0.05536953790374705,0.03922369633404027),
    @(0.9817236666666668,2.945171,0.04688329326954743,0.04832841473263862,12.68844067872533,114.195966108528,0.02893782949319845,0.03074927151949298),
    @(5.154927333333333,15.464782,0.2461792235003055,0.2537674037418691,66.62566248833066,599.6309623949759,0.1519494877090276,0.1614611785555975),
    @(1.8784265,3.756853,0.08970632314876403,0.06164760887349412,24.278016333784,145.668098002704,0.05536953790374705,0.03922369633404027),
    @(12.924656,38.773968,0.617231160081383,0.6362565726519981,167.046732718336,1503.420594465024,0.3809743049754098,0.4048224262428673)
)

$ws.Range("G2:J17").Value = (ConvertTo-Excel2DArray $gjJagged)
$ws.Range("M2:T17").Value = (ConvertTo-Excel2DArray $mtJagged)
